$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-06-17"

# Update the header label shared string for the current-year column
$ws.Range("I1").Value = "2022 (through 06-17)"

# Update June (row 7) value for 2022 column
$ws.Range("I7").Value = 82

# Update the Total row (row 14) value for 2022 column
$ws.Range("I14").Value = 745
